# Auto-generated edit script: updates market-data columns (H:N) for specific rows
# across multiple sheets, per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16748
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53736
$ws.Range("M72").ClearContents()

$ws.Range("H74").Value = 4859.3
$ws.Range("I74").Value = 4348.8335
$ws.Range("J74").Value = 5625
$ws.Range("K74").Value = 4348.8335
$ws.Range("L74").Value = 5625
$ws.Range("M74").Value = -3412.8335
$ws.Range("N74").Value = -7497

$ws.Range("H77").Value = 4859.3
$ws.Range("I77").Value = 4348.8335
$ws.Range("J77").Value = 5625
$ws.Range("K77").Value = 21744.1675
$ws.Range("L77").Value = 28125
$ws.Range("M77").Value = -17064.1675
$ws.Range("N77").Value = -37485

$ws.Range("H100").Value = 1869.1666
$ws.Range("I100").Value = 1903
$ws.Range("J100").Value = 1700
$ws.Range("K100").Value = 1903
$ws.Range("L100").Value = 1700
$ws.Range("M100").Value = -1362
$ws.Range("N100").Value = -2782

$ws.Range("H131").Value = 4311.4346
$ws.Range("I131").Value = 1253.6666
$ws.Range("J131").Value = 6277.143
$ws.Range("K131").Value = 3760.9998
$ws.Range("L131").Value = 18831.429
$ws.Range("M131").Value = 1279.0002
$ws.Range("N131").Value = -28911.429

$ws.Range("H137").Value = 4685.857
$ws.Range("I137").Value = 6087.75
$ws.Range("J137").Value = 2816.6667
$ws.Range("K137").Value = 18263.25
$ws.Range("L137").Value = 8450.000100000001
$ws.Range("M137").Value = -15713.25
$ws.Range("N137").Value = -13550.0001

$ws.Range("H138").Value = 29654.084
$ws.Range("I138").Value = 5239.6
$ws.Range("K138").Value = 15718.8
$ws.Range("M138").Value = -10578.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3663.3333
$ws.Range("I63").Value = 2995
$ws.Range("K63").Value = 2995
$ws.Range("M63").Value = -2309

$ws.Range("H66").Value = 3663.3333
$ws.Range("I66").Value = 2995
$ws.Range("K66").Value = 14975
$ws.Range("M66").Value = -11543

$ws.Range("H97").Value = 1145.5
$ws.Range("I97").Value = 1193.375
$ws.Range("J97").Value = 1049.75
$ws.Range("K97").Value = 1193.375
$ws.Range("L97").Value = 1049.75
$ws.Range("M97").Value = -697.375
$ws.Range("N97").Value = -2041.75

$ws.Range("H132").Value = 12953.818
$ws.Range("I132").Value = 27644.75
$ws.Range("J132").Value = 4559
$ws.Range("K132").Value = 82934.25
$ws.Range("L132").Value = 13677
$ws.Range("M132").Value = -80404.25
$ws.Range("N132").Value = -18737

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10665.5
$ws.Range("I31").Value = 11281
$ws.Range("K31").Value = 11281
$ws.Range("M31").Value = -10986

$ws.Range("H34").Value = 10665.5
$ws.Range("I34").Value = 11281
$ws.Range("K34").Value = 11281
$ws.Range("M34").Value = -11079

$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376

$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880

$ws.Range("H132").Value = 2429.7812
$ws.Range("I132").Value = 1827.2609
$ws.Range("J132").Value = 3969.5557
$ws.Range("K132").Value = 5481.7827
$ws.Range("L132").Value = 11908.6671
$ws.Range("M132").Value = -2951.7827
$ws.Range("N132").Value = -16968.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 11118123
$ws.Range("I5").Value = 338.91666
$ws.Range("K5").Value = 1016.74998
$ws.Range("M5").Value = -904.7499799999999

$ws.Range("H114").Value = 698.36365
$ws.Range("I114").Value = 363.14285
$ws.Range("J114").Value = 1285
$ws.Range("K114").Value = 1089.42855
$ws.Range("L114").Value = 3855
$ws.Range("M114").Value = 2164.57145
$ws.Range("N114").Value = -10363

$ws.Range("H122").Value = 1166.3334
$ws.Range("I122").Value = 238.2
$ws.Range("J122").Value = 1630.4
$ws.Range("K122").Value = 2143.8
$ws.Range("L122").Value = 14673.6
$ws.Range("M122").Value = 306.2000000000003
$ws.Range("N122").Value = -19573.6

$ws.Range("H135").Value = 11118123
$ws.Range("I135").Value = 338.91666
$ws.Range("K135").Value = 3050.24994
$ws.Range("M135").Value = -515.2499399999997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 19320
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 19320
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 19320
$ws.Range("N54").Value = -20608
$ws.Range("M54").ClearContents()

$ws.Range("H82").Value = 1678.6
$ws.Range("I82").Value = 1698
$ws.Range("J82").Value = 1601
$ws.Range("K82").Value = 1698
$ws.Range("L82").Value = 1601
$ws.Range("M82").Value = -1337
$ws.Range("N82").Value = -2323

$ws.Range("H85").Value = 1678.6
$ws.Range("I85").Value = 1698
$ws.Range("J85").Value = 1601
$ws.Range("K85").Value = 1698
$ws.Range("L85").Value = 1601
$ws.Range("M85").Value = -450
$ws.Range("N85").Value = -4097

$ws.Range("H122").Value = 6376.829
$ws.Range("I122").Value = 5613.3335
$ws.Range("J122").Value = 8459.091
$ws.Range("K122").Value = 16840.0005
$ws.Range("L122").Value = 25377.273
$ws.Range("M122").Value = -14390.0005
$ws.Range("N122").Value = -30277.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 15000
$ws.Range("I52").Value = 15000
$ws.Range("K52").Value = 15000
$ws.Range("M52").Value = -14774

$ws.Range("H107").Value = 2821.5293
$ws.Range("I107").Value = 769.25
$ws.Range("J107").Value = 4645.778
$ws.Range("K107").Value = 2307.75
$ws.Range("L107").Value = 13937.334
$ws.Range("M107").Value = -387.75
$ws.Range("N107").Value = -17777.334

$ws.Range("H112").Value = 59077.832
$ws.Range("J112").Value = 59077.832
$ws.Range("L112").Value = 59077.832
$ws.Range("N112").Value = -62031.832

$ws.Range("H122").Value = 4911.037
$ws.Range("I122").Value = 1311.625
$ws.Range("J122").Value = 10146.546
$ws.Range("K122").Value = 3934.875
$ws.Range("L122").Value = 30439.638
$ws.Range("M122").Value = -1484.875
$ws.Range("N122").Value = -35339.638

$ws.Range("H132").Value = 3070.4814
$ws.Range("I132").Value = 2873.9473
$ws.Range("J132").Value = 3537.25
$ws.Range("K132").Value = 8621.841899999999
$ws.Range("L132").Value = 10611.75
$ws.Range("M132").Value = -6091.841899999999
$ws.Range("N132").Value = -15671.75

